$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Pre-seed the shared-strings table in the exact order the new strings
#    appear in the target workbook. Shared strings are interned in
#    first-use order, and that order does not always match the final
#    top-to-bottom row order (some later rows reuse strings that were
#    first introduced by an earlier-indexed, but later-written, row).
#    Writing them once into a scratch area (far below the used range)
#    registers them in the shared-strings table with the right indices;
#    the scratch rows are removed again afterwards. Because the real
#    rows below still reference the same text, the strings remain in
#    the table (just no longer referenced by the scratch cells).
# ---------------------------------------------------------------------------
$ws.Cells.Item(2000, 1).Value2 = 'Clase 07'
$ws.Cells.Item(2001, 1).Value2 = 'Codigos de estado (1xx, 2xx, 3xx, 4xx, 5xx)'
$ws.Cells.Item(2002, 1).Value2 = 'Lo ideal es agregar un codigo de estado en la respuesta. Para que express no lo haga por defecto'
$ws.Cells.Item(2003, 1).Value2 = 'Formatos JSON y XML (Este ultimo casi no se utiliza en la actualidad)'
$ws.Cells.Item(2004, 1).Value2 = 'API REST'
$ws.Cells.Item(2005, 1).Value2 = 'Modelo de una API REST'
$ws.Cells.Item(2006, 1).Value2 = 'Metodos de peticion'
$ws.Cells.Item(2007, 1).Value2 = 'Copiando el user.manager con los metodos getUsers y createUser'
$ws.Cells.Item(2008, 1).Value2 = 'importar uuid'
$ws.Cells.Item(2009, 1).Value2 = 'modificacion createUser'
$ws.Cells.Item(2010, 1).Value2 = 'creacion de getUserById'
$ws.Cells.Item(2011, 1).Value2 = 'creacion de updateUser'
$ws.Cells.Item(2012, 1).Value2 = 'creacion de deleteUser'
$ws.Cells.Item(2013, 1).Value2 = 'Middleware para datos que se envian por parametros, por URL'
$ws.Cells.Item(2014, 1).Value2 = 'Endpoint app.get(''/users'''
$ws.Cells.Item(2015, 1).Value2 = 'Endpoint app.get(''/users/:id'''
$ws.Cells.Item(2016, 1).Value2 = 'Endpoint app.post(''/users'''
$ws.Cells.Item(2017, 1).Value2 = 'Explicacion postman'
$ws.Cells.Item(2018, 1).Value2 = 'Endpoint app.put('''
$ws.Cells.Item(2019, 1).Value2 = 'Endpoint app.delete('''
$ws.Cells.Item(2020, 1).Value2 = 'Clase 08'
$ws.Cells.Item(2021, 1).Value2 = 'Router en Express'
$ws.Cells.Item(2022, 1).Value2 = 'Pasando las Rutas del server.js al user.router.js'
$ws.Cells.Item(2023, 1).Value2 = 'Revisando con postman despues de los cambios'
$ws.Cells.Item(2024, 1).Value2 = 'concepto modularizacion'
$ws.Cells.Item(2025, 1).Value2 = 'Hands on lab "Users y Pets"'
$ws.Cells.Item(2026, 1).Value2 = 'Cada Router deberia tener su propio manager. El ejemplo de Pets se hizo asi solo por explicacion'
$ws.Cells.Item(2027, 1).Value2 = 'Middleware static  app.use(express.static(''public''))'
$ws.Cells.Item(2028, 1).Value2 = 'Hay Middlewares a nivel de aplicación, a nivel de enrutador y a nivel de endpoint'
$ws.Cells.Item(2029, 1).Value2 = 'Ejemplo de como funciona un Middleware'
$ws.Cells.Item(2030, 1).Value2 = '__dirname'
$ws.Cells.Item(2031, 1).Value2 = 'Archivos estaticos con express.static (imagen, html) '
$ws.Cells.Item(2032, 1).Value2 = 'definicion ¿Qué es un middleware?'
$ws.Cells.Item(2033, 1).Value2 = 'Multer es un Middleware que permite subir archivos al servidor'
$ws.Cells.Item(2034, 1).Value2 = 'Tipos de Middleware'
$ws.Cells.Item(2035, 1).Value2 = 'Multer - Es un middleware de terceros'
$ws.Cells.Item(2036, 1).Value2 = 'Instalacion de multer'
$ws.Cells.Item(2037, 1).Value2 = 'Archivo multer.js'
$ws.Cells.Item(2038, 1).Value2 = 'Subiendo una imagen como adjunto con postman'

# ---------------------------------------------------------------------------
# 2) Append the "Clase 07" / "Clase 08" sections (rows 113-152), copying
#    cell formatting from the existing header/data template rows.
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

function Add-HeaderRow([int]$row, [string]$text) {
    [void]$ws.Range("A8:B8").Copy()
    [void]$ws.Range("A$row`:B$row").PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($row, 1).Value2 = $text
}

function Add-DataRow([int]$row, [double]$time, [string]$text) {
    [void]$ws.Range("A9:B9").Copy()
    [void]$ws.Range("A$row`:B$row").PasteSpecial($xlPasteFormats)
    $ws.Cells.Item($row, 1).Value2 = $time
    $ws.Cells.Item($row, 2).Value2 = $text
}

Add-HeaderRow 113 'Clase 07'
Add-DataRow 114 0.0027777777777777779 'Codigos de estado (1xx, 2xx, 3xx, 4xx, 5xx)'
Add-DataRow 115 0.0048611111111111112 'Lo ideal es agregar un codigo de estado en la respuesta. Para que express no lo haga por defecto'
Add-DataRow 116 0.0072916666666666668 'API REST'
Add-DataRow 117 0.010069444444444445 'Formatos JSON y XML (Este ultimo casi no se utiliza en la actualidad)'
Add-DataRow 118 0.011111111111111112 'Modelo de una API REST'
Add-DataRow 119 0.012500000000000001 'Metodos de peticion'
Add-DataRow 120 0.014930555555555556 'Copiando el user.manager con los metodos getUsers y createUser'
Add-DataRow 121 0.016319444444444445 'importar uuid'
Add-DataRow 122 0.016666666666666666 'modificacion createUser'
Add-DataRow 123 0.02013888888888889 'creacion de getUserById'
Add-DataRow 124 0.021874999999999999 'creacion de updateUser'
Add-DataRow 125 0.025347222222222222 'creacion de deleteUser'
Add-DataRow 126 0.027430555555555555 'Middleware para datos que se envian por parametros, por URL'
Add-DataRow 127 0.028125000000000001 'Endpoint app.get(''/users'''
Add-DataRow 128 0.030555555555555555 'Endpoint app.get(''/users/:id'''
Add-DataRow 129 0.036111111111111108 'Endpoint app.post(''/users'''
Add-DataRow 130 0.044791666666666667 'Explicacion postman'
Add-DataRow 131 0.057291666666666664 'Endpoint app.put('''
Add-DataRow 132 0.061805555555555558 'Endpoint app.delete('''
Add-HeaderRow 133 'Clase 08'
Add-DataRow 134 0.0013888888888888889 'Router en Express'
Add-DataRow 135 0.0059027777777777776 'Pasando las Rutas del server.js al user.router.js'
Add-DataRow 136 0.012500000000000001 'Revisando con postman despues de los cambios'
Add-DataRow 137 0.012847222222222222 'concepto modularizacion'
Add-DataRow 138 0.013541666666666667 'Hands on lab "Users y Pets"'
Add-DataRow 139 0.02013888888888889 'Cada Router deberia tener su propio manager. El ejemplo de Pets se hizo asi solo por explicacion'
Add-DataRow 140 0.020486111111111111 'Archivos estaticos con express.static (imagen, html) '
Add-DataRow 141 0.020833333333333332 'Middleware static  app.use(express.static(''public''))'
Add-DataRow 142 0.021527777777777778 'Hay Middlewares a nivel de aplicación, a nivel de enrutador y a nivel de endpoint'
Add-DataRow 143 0.022916666666666665 'Ejemplo de como funciona un Middleware'
Add-DataRow 144 0.03125 'Middleware static  app.use(express.static(''public''))'
Add-DataRow 145 0.034722222222222224 '__dirname'
Add-DataRow 146 0.049652777777777775 'definicion ¿Qué es un middleware?'
Add-DataRow 147 0.051388888888888887 'Multer es un Middleware que permite subir archivos al servidor'
Add-DataRow 148 0.052083333333333336 'Tipos de Middleware'
Add-DataRow 149 0.053124999999999999 'Multer - Es un middleware de terceros'
Add-DataRow 150 0.054166666666666669 'Instalacion de multer'
Add-DataRow 151 0.055208333333333331 'Archivo multer.js'
Add-DataRow 152 0.062847222222222221 'Subiendo una imagen como adjunto con postman'

# ---------------------------------------------------------------------------
# 3) Remove the temporary scratch rows used purely for seeding order.
# ---------------------------------------------------------------------------
[void]$ws.Range("A2000:A2038").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 4) Mirror the final selection recorded in the workbook. (The scroll /
#    topLeftCell offset itself is not persisted by this engine unless a
#    freeze/split pane is in use, so only the active-cell selection is
#    reproduced here.)
# ---------------------------------------------------------------------------
[void]$ws.Range("B152").Select()

Write-Host "Added Clase 07 and Clase 08 rows (113-152)."
